$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 6946332
$ws.Cells.Item(17, 10).Value = 7577629.5
$ws.Cells.Item(17, 12).Value = 22732888.5
$ws.Cells.Item(17, 14).Value = -22733224.5
$ws.Cells.Item(62, 8).Value = 2736.65
$ws.Cells.Item(62, 9).Value = 2718.889
$ws.Cells.Item(62, 10).Value = 2896.5
$ws.Cells.Item(62, 11).Value = 2718.889
$ws.Cells.Item(62, 12).Value = 2896.5
$ws.Cells.Item(62, 13).Value = -2094.889
$ws.Cells.Item(62, 14).Value = -4144.5
$ws.Cells.Item(64, 8).Value = 166674670
$ws.Cells.Item(64, 9).Value = 9334
$ws.Cells.Item(64, 10).Value = 333340000
$ws.Cells.Item(64, 11).Value = 9334
$ws.Cells.Item(64, 12).Value = 333340000
$ws.Cells.Item(64, 13).Value = -9086
$ws.Cells.Item(64, 14).Value = -333340496
$ws.Cells.Item(65, 8).Value = 2736.65
$ws.Cells.Item(65, 9).Value = 2718.889
$ws.Cells.Item(65, 10).Value = 2896.5
$ws.Cells.Item(65, 11).Value = 13594.445
$ws.Cells.Item(65, 12).Value = 14482.5
$ws.Cells.Item(65, 13).Value = -10474.445
$ws.Cells.Item(65, 14).Value = -20722.5
$ws.Cells.Item(67, 8).Value = 166674670
$ws.Cells.Item(67, 9).Value = 9334
$ws.Cells.Item(67, 10).Value = 333340000
$ws.Cells.Item(67, 11).Value = 9334
$ws.Cells.Item(67, 12).Value = 333340000
$ws.Cells.Item(67, 13).Value = -8476
$ws.Cells.Item(67, 14).Value = -333341716
$ws.Cells.Item(74, 8).Value = 17517.53
$ws.Cells.Item(74, 9).Value = 17842.857
$ws.Cells.Item(74, 11).Value = 17842.857
$ws.Cells.Item(74, 13).Value = -16906.857
$ws.Cells.Item(77, 8).Value = 17517.53
$ws.Cells.Item(77, 9).Value = 17842.857
$ws.Cells.Item(77, 11).Value = 89214.285
$ws.Cells.Item(77, 13).Value = -84534.285
$ws.Cells.Item(86, 8).Value = 111113370
$ws.Cells.Item(86, 9).Value = 142858480
$ws.Cells.Item(86, 11).Value = 142858480
$ws.Cells.Item(86, 13).Value = -142857357
$ws.Cells.Item(88, 8).Value = 11111989
$ws.Cells.Item(88, 9).Value = 33333966
$ws.Cells.Item(88, 10).Value = 1000
$ws.Cells.Item(88, 11).Value = 33333966
$ws.Cells.Item(88, 12).Value = 1000
$ws.Cells.Item(88, 13).Value = -33333560
$ws.Cells.Item(88, 14).Value = -1812
$ws.Cells.Item(89, 8).Value = 111113370
$ws.Cells.Item(89, 9).Value = 142858480
$ws.Cells.Item(89, 11).Value = 714292400
$ws.Cells.Item(89, 13).Value = -714286784
$ws.Cells.Item(91, 8).Value = 11111989
$ws.Cells.Item(91, 9).Value = 33333966
$ws.Cells.Item(91, 10).Value = 1000
$ws.Cells.Item(91, 11).Value = 33333966
$ws.Cells.Item(91, 12).Value = 1000
$ws.Cells.Item(91, 13).Value = -33332562
$ws.Cells.Item(91, 14).Value = -3808
$ws.Cells.Item(92, 8).Value = 559.53845
$ws.Cells.Item(92, 9).Value = 559.53845
$ws.Cells.Item(92, 11).Value = 559.53845
$ws.Cells.Item(92, 13).Value = 688.46155
$ws.Cells.Item(96, 8).Value = 791.9231
$ws.Cells.Item(96, 9).Value = 398.75
$ws.Cells.Item(96, 10).Value = 966.6667
$ws.Cells.Item(96, 11).Value = 1196.25
$ws.Cells.Item(96, 12).Value = 2900.0001
$ws.Cells.Item(96, 13).Value = 176.75
$ws.Cells.Item(96, 14).Value = -5646.0001
$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 10).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 12).Value = 0
$ws.Cells.Item(100, 13).ClearContents()
$ws.Cells.Item(100, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 404.3684
$ws.Cells.Item(107, 9).Value = 495.15384
$ws.Cells.Item(107, 11).Value = 495.15384
$ws.Cells.Item(107, 13).Value = 1424.84616
$ws.Cells.Item(108, 8).Value = 73135.2
$ws.Cells.Item(108, 10).Value = 73135.2
$ws.Cells.Item(108, 12).Value = 73135.2
$ws.Cells.Item(108, 14).Value = -80815.2
$ws.Cells.Item(116, 8).Value = 9594.208
$ws.Cells.Item(116, 9).Value = 2831
$ws.Cells.Item(116, 10).Value = 10000
$ws.Cells.Item(116, 11).Value = 2831
$ws.Cells.Item(116, 12).Value = 10000
$ws.Cells.Item(116, 13).Value = 611
$ws.Cells.Item(116, 14).Value = -16884
$ws.Cells.Item(125, 8).Value = 1215
$ws.Cells.Item(125, 10).Value = 1215
$ws.Cells.Item(125, 12).Value = 10935
$ws.Cells.Item(125, 14).Value = -15855
$ws.Cells.Item(135, 8).Value = 764.6111
$ws.Cells.Item(135, 9).Value = 295.5
$ws.Cells.Item(135, 11).Value = 2659.5
$ws.Cells.Item(135, 13).Value = -124.5
$ws.Cells.Item(137, 8).Value = 5354.76
$ws.Cells.Item(137, 9).Value = 1568.2858
$ws.Cells.Item(137, 10).Value = 10173.909
$ws.Cells.Item(137, 11).Value = 4704.857400000001
$ws.Cells.Item(137, 12).Value = 30521.727
$ws.Cells.Item(137, 13).Value = -2154.857400000001
$ws.Cells.Item(137, 14).Value = -35621.727
$ws.Cells.Item(138, 8).Value = 332375.06
$ws.Cells.Item(138, 10).Value = 513498.78
$ws.Cells.Item(138, 12).Value = 1540496.34
$ws.Cells.Item(138, 14).Value = -1550776.34
$ws.Cells.Item(141, 8).Value = 5113.5
$ws.Cells.Item(141, 9).Value = 5116.3335
$ws.Cells.Item(141, 10).Value = 5105
$ws.Cells.Item(141, 11).Value = 15349.0005
$ws.Cells.Item(141, 12).Value = 15315
$ws.Cells.Item(141, 13).Value = -10169.0005
$ws.Cells.Item(141, 14).Value = -25675

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4193.608
$ws.Cells.Item(32, 9).Value = 4257.48
$ws.Cells.Item(32, 10).Value = 1000
$ws.Cells.Item(32, 11).Value = 4257.48
$ws.Cells.Item(32, 12).Value = 1000
$ws.Cells.Item(32, 13).Value = -3970.48
$ws.Cells.Item(32, 14).Value = -1574
$ws.Cells.Item(45, 8).Value = 17203.908
$ws.Cells.Item(45, 9).Value = 22382.87
$ws.Cells.Item(45, 10).Value = 5292.3
$ws.Cells.Item(45, 11).Value = 22382.87
$ws.Cells.Item(45, 12).Value = 5292.3
$ws.Cells.Item(45, 13).Value = -22005.87
$ws.Cells.Item(45, 14).Value = -6046.3
$ws.Cells.Item(61, 8).Value = 3951.3225
$ws.Cells.Item(61, 9).Value = 2259.8262
$ws.Cells.Item(61, 11).Value = 2259.8262
$ws.Cells.Item(61, 13).Value = -2047.8262
$ws.Cells.Item(88, 8).Value = 5879.778
$ws.Cells.Item(88, 9).Value = 4246.6
$ws.Cells.Item(88, 10).Value = 7921.25
$ws.Cells.Item(88, 11).Value = 4246.6
$ws.Cells.Item(88, 12).Value = 7921.25
$ws.Cells.Item(88, 13).Value = -3840.6
$ws.Cells.Item(88, 14).Value = -8733.25
$ws.Cells.Item(91, 8).Value = 5879.778
$ws.Cells.Item(91, 9).Value = 4246.6
$ws.Cells.Item(91, 10).Value = 7921.25
$ws.Cells.Item(91, 11).Value = 4246.6
$ws.Cells.Item(91, 12).Value = 7921.25
$ws.Cells.Item(91, 13).Value = -2842.6
$ws.Cells.Item(91, 14).Value = -10729.25
$ws.Cells.Item(97, 8).Value = 1359.2142
$ws.Cells.Item(97, 9).Value = 1359.2142
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 11).Value = 1359.2142
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 13).Value = -863.2141999999999
$ws.Cells.Item(97, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 5203.5654
$ws.Cells.Item(102, 9).Value = 5341.933
$ws.Cells.Item(102, 11).Value = 5341.933
$ws.Cells.Item(102, 13).Value = -3719.933
$ws.Cells.Item(122, 8).Value = 4269
$ws.Cells.Item(122, 9).Value = 3684.1904
$ws.Cells.Item(122, 10).Value = 6725.2
$ws.Cells.Item(122, 11).Value = 11052.5712
$ws.Cells.Item(122, 12).Value = 20175.6
$ws.Cells.Item(122, 13).Value = -8602.5712
$ws.Cells.Item(122, 14).Value = -25075.6
$ws.Cells.Item(132, 8).Value = 2537.9143
$ws.Cells.Item(132, 9).Value = 1886.4073
$ws.Cells.Item(132, 10).Value = 4736.75
$ws.Cells.Item(132, 11).Value = 5659.2219
$ws.Cells.Item(132, 12).Value = 14210.25
$ws.Cells.Item(132, 13).Value = -3129.2219
$ws.Cells.Item(132, 14).Value = -19270.25
$ws.Cells.Item(136, 8).Value = 3951.3225
$ws.Cells.Item(136, 9).Value = 2259.8262
$ws.Cells.Item(136, 11).Value = 6779.4786
$ws.Cells.Item(136, 13).Value = -4229.4786

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 25010008
$ws.Cells.Item(20, 9).Value = 31261136
$ws.Cells.Item(20, 10).Value = 5499.5
$ws.Cells.Item(20, 11).Value = 31261136
$ws.Cells.Item(20, 12).Value = 5499.5
$ws.Cells.Item(20, 13).Value = -31260889
$ws.Cells.Item(20, 14).Value = -5993.5
$ws.Cells.Item(86, 8).Value = 3745.8462
$ws.Cells.Item(86, 9).Value = 3245.182
$ws.Cells.Item(86, 11).Value = 3245.182
$ws.Cells.Item(86, 13).Value = -2122.182
$ws.Cells.Item(89, 8).Value = 3745.8462
$ws.Cells.Item(89, 9).Value = 3245.182
$ws.Cells.Item(89, 11).Value = 16225.91
$ws.Cells.Item(89, 13).Value = -10609.91
$ws.Cells.Item(94, 8).Value = 50000176
$ws.Cells.Item(94, 9).Value = 55555744
$ws.Cells.Item(94, 10).Value = 74
$ws.Cells.Item(94, 11).Value = 55555744
$ws.Cells.Item(94, 12).Value = 74
$ws.Cells.Item(94, 13).Value = -55555293
$ws.Cells.Item(94, 14).Value = -976
$ws.Cells.Item(105, 8).Value = 11307052
$ws.Cells.Item(105, 9).Value = 669204.06
$ws.Cells.Item(105, 10).Value = 31253018
$ws.Cells.Item(105, 11).Value = 669204.06
$ws.Cells.Item(105, 12).Value = 31253018
$ws.Cells.Item(105, 13).Value = -667457.06
$ws.Cells.Item(105, 14).Value = -31256512
$ws.Cells.Item(107, 8).Value = 2405171.5
$ws.Cells.Item(107, 9).Value = 2959710
$ws.Cells.Item(107, 11).Value = 2959710
$ws.Cells.Item(107, 13).Value = -2957790
$ws.Cells.Item(134, 8).Value = 1913.0834
$ws.Cells.Item(134, 9).Value = 1312.3667
$ws.Cells.Item(134, 11).Value = 3937.1001
$ws.Cells.Item(134, 13).Value = -1402.1001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1807.3429
$ws.Cells.Item(16, 9).Value = 1748.8
$ws.Cells.Item(16, 11).Value = 1748.8
$ws.Cells.Item(16, 13).Value = -1461.8
$ws.Cells.Item(31, 8).Value = 5297
$ws.Cells.Item(31, 9).Value = 4898.9165
$ws.Cells.Item(31, 10).Value = 5827.778
$ws.Cells.Item(31, 11).Value = 4898.9165
$ws.Cells.Item(31, 12).Value = 5827.778
$ws.Cells.Item(31, 13).Value = -4603.9165
$ws.Cells.Item(31, 14).Value = -6417.778
$ws.Cells.Item(34, 8).Value = 5297
$ws.Cells.Item(34, 9).Value = 4898.9165
$ws.Cells.Item(34, 10).Value = 5827.778
$ws.Cells.Item(34, 11).Value = 4898.9165
$ws.Cells.Item(34, 12).Value = 5827.778
$ws.Cells.Item(34, 13).Value = -4696.9165
$ws.Cells.Item(34, 14).Value = -6231.778
$ws.Cells.Item(86, 8).Value = 1013144.1
$ws.Cells.Item(86, 9).Value = 3058.6667
$ws.Cells.Item(86, 10).Value = 2907054.5
$ws.Cells.Item(86, 11).Value = 3058.6667
$ws.Cells.Item(86, 12).Value = 2907054.5
$ws.Cells.Item(86, 13).Value = -1935.6667
$ws.Cells.Item(86, 14).Value = -2909300.5
$ws.Cells.Item(89, 8).Value = 1013144.1
$ws.Cells.Item(89, 9).Value = 3058.6667
$ws.Cells.Item(89, 10).Value = 2907054.5
$ws.Cells.Item(89, 11).Value = 15293.3335
$ws.Cells.Item(89, 12).Value = 14535272.5
$ws.Cells.Item(89, 13).Value = -9677.3335
$ws.Cells.Item(89, 14).Value = -14546504.5
$ws.Cells.Item(107, 8).Value = 2174333
$ws.Cells.Item(107, 9).Value = 2500408
$ws.Cells.Item(107, 11).Value = 2500408
$ws.Cells.Item(107, 13).Value = -2498488
$ws.Cells.Item(113, 8).Value = 1807.3429
$ws.Cells.Item(113, 9).Value = 1748.8
$ws.Cells.Item(113, 11).Value = 1748.8
$ws.Cells.Item(113, 13).Value = 421.2
$ws.Cells.Item(132, 8).Value = 17244742
$ws.Cells.Item(132, 9).Value = 23812198
$ws.Cells.Item(132, 10).Value = 5170.375
$ws.Cells.Item(132, 11).Value = 71436594
$ws.Cells.Item(132, 12).Value = 15511.125
$ws.Cells.Item(132, 13).Value = -71434064
$ws.Cells.Item(132, 14).Value = -20571.125
$ws.Cells.Item(134, 8).Value = 4661.1787
$ws.Cells.Item(134, 9).Value = 4356.478
$ws.Cells.Item(134, 11).Value = 13069.434
$ws.Cells.Item(134, 13).Value = -10534.434

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 7112.5
$ws.Cells.Item(62, 10).Value = 7112.5
$ws.Cells.Item(62, 12).Value = 21337.5
$ws.Cells.Item(62, 14).Value = -22709.5
$ws.Cells.Item(63, 8).Value = 3650
$ws.Cells.Item(63, 10).Value = 5500
$ws.Cells.Item(63, 12).Value = 16500
$ws.Cells.Item(63, 14).Value = -17998
$ws.Cells.Item(65, 8).Value = 7112.5
$ws.Cells.Item(65, 10).Value = 7112.5
$ws.Cells.Item(65, 12).Value = 64012.5
$ws.Cells.Item(65, 14).Value = -70876.5
$ws.Cells.Item(66, 8).Value = 3650
$ws.Cells.Item(66, 10).Value = 5500
$ws.Cells.Item(66, 12).Value = 49500
$ws.Cells.Item(66, 14).Value = -56988
$ws.Cells.Item(80, 8).Value = 3906
$ws.Cells.Item(80, 10).Value = 3906
$ws.Cells.Item(80, 12).Value = 11718
$ws.Cells.Item(80, 14).Value = -13590
$ws.Cells.Item(83, 8).Value = 3906
$ws.Cells.Item(83, 10).Value = 3906
$ws.Cells.Item(83, 12).Value = 35154
$ws.Cells.Item(83, 14).Value = -44514
$ws.Cells.Item(103, 8).Value = 2278
$ws.Cells.Item(103, 9).Value = 547
$ws.Cells.Item(103, 10).Value = 4009
$ws.Cells.Item(103, 11).Value = 1641
$ws.Cells.Item(103, 12).Value = 12027
$ws.Cells.Item(103, 13).Value = -762
$ws.Cells.Item(103, 14).Value = -13785
$ws.Cells.Item(114, 8).Value = 990.6667
$ws.Cells.Item(114, 9).Value = 749.7273
$ws.Cells.Item(114, 10).Value = 1369.2858
$ws.Cells.Item(114, 11).Value = 2249.1819
$ws.Cells.Item(114, 12).Value = 4107.857400000001
$ws.Cells.Item(114, 13).Value = 1004.8181
$ws.Cells.Item(114, 14).Value = -10615.8574
$ws.Cells.Item(115, 8).Value = 6987.4443
$ws.Cells.Item(115, 9).Value = 1944
$ws.Cells.Item(115, 10).Value = 8428.429
$ws.Cells.Item(115, 11).Value = 5832
$ws.Cells.Item(115, 12).Value = 25285.287
$ws.Cells.Item(115, 13).Value = -4657
$ws.Cells.Item(115, 14).Value = -27635.287
$ws.Cells.Item(122, 8).Value = 1109
$ws.Cells.Item(122, 9).Value = 1072.2
$ws.Cells.Item(122, 10).Value = 1132
$ws.Cells.Item(122, 11).Value = 9649.800000000001
$ws.Cells.Item(122, 12).Value = 10188
$ws.Cells.Item(122, 13).Value = -7199.800000000001
$ws.Cells.Item(122, 14).Value = -15088
$ws.Cells.Item(126, 8).Value = 3049.5
$ws.Cells.Item(126, 9).Value = 3049.5
$ws.Cells.Item(126, 11).Value = 9148.5
$ws.Cells.Item(126, 13).Value = -4208.5
$ws.Cells.Item(131, 8).Value = 2078.5217
$ws.Cells.Item(131, 10).Value = 2259.2727
$ws.Cells.Item(131, 12).Value = 6777.8181
$ws.Cells.Item(131, 14).Value = -16857.8181
$ws.Cells.Item(132, 8).Value = 2860.6538
$ws.Cells.Item(132, 10).Value = 3471.2942
$ws.Cells.Item(132, 12).Value = 31241.6478
$ws.Cells.Item(132, 14).Value = -36301.6478

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(27, 8).Value = 1337500
$ws.Cells.Item(27, 9).Value = 5050000
$ws.Cells.Item(27, 11).Value = 5050000
$ws.Cells.Item(27, 13).Value = -5049834
$ws.Cells.Item(70, 8).Value = 87999.54
$ws.Cells.Item(70, 9).Value = 147517.64
$ws.Cells.Item(70, 10).Value = 4674.2
$ws.Cells.Item(70, 11).Value = 147517.64
$ws.Cells.Item(70, 12).Value = 4674.2
$ws.Cells.Item(70, 13).Value = -147247.64
$ws.Cells.Item(70, 14).Value = -5214.2
$ws.Cells.Item(73, 8).Value = 87999.54
$ws.Cells.Item(73, 9).Value = 147517.64
$ws.Cells.Item(73, 10).Value = 4674.2
$ws.Cells.Item(73, 11).Value = 147517.64
$ws.Cells.Item(73, 12).Value = 4674.2
$ws.Cells.Item(73, 13).Value = -146581.64
$ws.Cells.Item(73, 14).Value = -6546.2
$ws.Cells.Item(80, 8).Value = 142867070
$ws.Cells.Item(80, 9).Value = 200009000
$ws.Cells.Item(80, 10).Value = 12250
$ws.Cells.Item(80, 11).Value = 200009000
$ws.Cells.Item(80, 12).Value = 12250
$ws.Cells.Item(80, 13).Value = -200008002
$ws.Cells.Item(80, 14).Value = -14246
$ws.Cells.Item(83, 8).Value = 142867070
$ws.Cells.Item(83, 9).Value = 200009000
$ws.Cells.Item(83, 10).Value = 12250
$ws.Cells.Item(83, 11).Value = 1000045000
$ws.Cells.Item(83, 12).Value = 61250
$ws.Cells.Item(83, 13).Value = -1000040008
$ws.Cells.Item(83, 14).Value = -71234
$ws.Cells.Item(102, 8).Value = 5222.1665
$ws.Cells.Item(102, 9).Value = 791.6
$ws.Cells.Item(102, 10).Value = 27375
$ws.Cells.Item(102, 11).Value = 791.6
$ws.Cells.Item(102, 12).Value = 27375
$ws.Cells.Item(102, 13).Value = 830.4
$ws.Cells.Item(102, 14).Value = -30619
$ws.Cells.Item(113, 8).Value = 3665.9092
$ws.Cells.Item(113, 9).Value = 3505.8276
$ws.Cells.Item(113, 11).Value = 3505.8276
$ws.Cells.Item(113, 13).Value = -1335.8276
$ws.Cells.Item(122, 8).Value = 10771.667
$ws.Cells.Item(122, 9).Value = 8982.667
$ws.Cells.Item(122, 11).Value = 26948.001
$ws.Cells.Item(122, 13).Value = -24498.001
$ws.Cells.Item(126, 8).Value = 7429.2383
$ws.Cells.Item(126, 9).Value = 4934.533
$ws.Cells.Item(126, 10).Value = 13666
$ws.Cells.Item(126, 11).Value = 14803.599
$ws.Cells.Item(126, 12).Value = 40998
$ws.Cells.Item(126, 13).Value = -12333.599
$ws.Cells.Item(126, 14).Value = -45938
$ws.Cells.Item(132, 8).Value = 1676.0869
$ws.Cells.Item(132, 9).Value = 1226.4117
$ws.Cells.Item(132, 10).Value = 2950.1667
$ws.Cells.Item(132, 11).Value = 3679.2351
$ws.Cells.Item(132, 12).Value = 8850.500100000001
$ws.Cells.Item(132, 13).Value = -1149.2351
$ws.Cells.Item(132, 14).Value = -13910.5001
$ws.Cells.Item(135, 8).Value = 68423.43
$ws.Cells.Item(135, 10).Value = 68423.43
$ws.Cells.Item(135, 12).Value = 68423.43
$ws.Cells.Item(135, 14).Value = -78563.43
$ws.Cells.Item(136, 8).Value = 9725.5
$ws.Cells.Item(136, 10).Value = 9725.5
$ws.Cells.Item(136, 12).Value = 29176.5
$ws.Cells.Item(136, 14).Value = -34276.5
$ws.Cells.Item(137, 8).Value = 80000
$ws.Cells.Item(137, 10).Value = 80000
$ws.Cells.Item(137, 12).Value = 80000
$ws.Cells.Item(137, 14).Value = -90200

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4195.5
$ws.Cells.Item(7, 9).Value = 3521
$ws.Cells.Item(7, 10).Value = 6443.8335
$ws.Cells.Item(7, 11).Value = 3521
$ws.Cells.Item(7, 12).Value = 6443.8335
$ws.Cells.Item(7, 13).Value = -3409
$ws.Cells.Item(7, 14).Value = -6667.8335
$ws.Cells.Item(16, 8).Value = 1070.3572
$ws.Cells.Item(16, 9).Value = 1090.5
$ws.Cells.Item(16, 10).Value = 949.5
$ws.Cells.Item(16, 11).Value = 1090.5
$ws.Cells.Item(16, 12).Value = 949.5
$ws.Cells.Item(16, 13).Value = -920.5
$ws.Cells.Item(16, 14).Value = -1289.5
$ws.Cells.Item(38, 8).Value = 30000
$ws.Cells.Item(38, 10).Value = 30000
$ws.Cells.Item(38, 12).Value = 30000
$ws.Cells.Item(38, 14).Value = -30820
$ws.Cells.Item(40, 8).Value = 47583.043
$ws.Cells.Item(40, 9).Value = 49217.305
$ws.Cells.Item(40, 10).Value = 9995
$ws.Cells.Item(40, 11).Value = 49217.305
$ws.Cells.Item(40, 12).Value = 9995
$ws.Cells.Item(40, 13).Value = -49081.305
$ws.Cells.Item(40, 14).Value = -10267
$ws.Cells.Item(61, 8).Value = 978.64105
$ws.Cells.Item(61, 9).Value = 882.05554
$ws.Cells.Item(61, 10).Value = 2137.6667
$ws.Cells.Item(61, 11).Value = 882.05554
$ws.Cells.Item(61, 12).Value = 2137.6667
$ws.Cells.Item(61, 13).Value = -680.05554
$ws.Cells.Item(61, 14).Value = -2541.6667
$ws.Cells.Item(82, 8).Value = 1575
$ws.Cells.Item(82, 9).Value = 1575
$ws.Cells.Item(82, 11).Value = 1575
$ws.Cells.Item(82, 13).Value = -1214
$ws.Cells.Item(85, 8).Value = 1575
$ws.Cells.Item(85, 9).Value = 1575
$ws.Cells.Item(85, 11).Value = 1575
$ws.Cells.Item(85, 13).Value = -327
$ws.Cells.Item(93, 8).Value = 1599.2858
$ws.Cells.Item(93, 9).Value = 1774.25
$ws.Cells.Item(93, 10).Value = 1366
$ws.Cells.Item(93, 11).Value = 1774.25
$ws.Cells.Item(93, 12).Value = 1366
$ws.Cells.Item(93, 13).Value = -526.25
$ws.Cells.Item(93, 14).Value = -3862
$ws.Cells.Item(113, 8).Value = 978.64105
$ws.Cells.Item(113, 9).Value = 882.05554
$ws.Cells.Item(113, 10).Value = 2137.6667
$ws.Cells.Item(113, 11).Value = 882.05554
$ws.Cells.Item(113, 12).Value = 2137.6667
$ws.Cells.Item(113, 13).Value = 1287.94446
$ws.Cells.Item(113, 14).Value = -6477.6667
$ws.Cells.Item(122, 8).Value = 3233.1667
$ws.Cells.Item(122, 9).Value = 4025.125
$ws.Cells.Item(122, 10).Value = 1649.25
$ws.Cells.Item(122, 11).Value = 12075.375
$ws.Cells.Item(122, 12).Value = 4947.75
$ws.Cells.Item(122, 13).Value = -9625.375
$ws.Cells.Item(122, 14).Value = -9847.75
$ws.Cells.Item(126, 8).Value = 4195.5
$ws.Cells.Item(126, 9).Value = 3521
$ws.Cells.Item(126, 10).Value = 6443.8335
$ws.Cells.Item(126, 11).Value = 10563
$ws.Cells.Item(126, 12).Value = 19331.5005
$ws.Cells.Item(126, 13).Value = -8093
$ws.Cells.Item(126, 14).Value = -24271.5005
$ws.Cells.Item(132, 8).Value = 3505.9033
$ws.Cells.Item(132, 9).Value = 2676.3655
$ws.Cells.Item(132, 10).Value = 7819.5
$ws.Cells.Item(132, 11).Value = 8029.0965
$ws.Cells.Item(132, 12).Value = 23458.5
$ws.Cells.Item(132, 13).Value = -5499.0965
$ws.Cells.Item(132, 14).Value = -28518.5
$ws.Cells.Item(136, 8).Value = 3717.3
$ws.Cells.Item(136, 9).Value = 3391.7144
$ws.Cells.Item(136, 10).Value = 4477
$ws.Cells.Item(136, 11).Value = 10175.1432
$ws.Cells.Item(136, 12).Value = 13431
$ws.Cells.Item(136, 13).Value = -7625.143199999999
$ws.Cells.Item(136, 14).Value = -18531

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 12500
$ws.Cells.Item(18, 10).Value = 12500
$ws.Cells.Item(18, 12).Value = 12500
$ws.Cells.Item(18, 14).Value = -12846
$ws.Cells.Item(81, 8).Value = 13000
$ws.Cells.Item(81, 9).Value = 20000
$ws.Cells.Item(81, 10).Value = 6000
$ws.Cells.Item(81, 11).Value = 40000
$ws.Cells.Item(81, 12).Value = 12000
$ws.Cells.Item(81, 13).Value = -38939
$ws.Cells.Item(81, 14).Value = -14122
$ws.Cells.Item(84, 8).Value = 13000
$ws.Cells.Item(84, 9).Value = 20000
$ws.Cells.Item(84, 10).Value = 6000
$ws.Cells.Item(84, 11).Value = 200000
$ws.Cells.Item(84, 12).Value = 60000
$ws.Cells.Item(84, 13).Value = -194696
$ws.Cells.Item(84, 14).Value = -70608
$ws.Cells.Item(96, 8).Value = 5907
$ws.Cells.Item(96, 9).Value = 5929.1665
$ws.Cells.Item(96, 10).Value = 5869
$ws.Cells.Item(96, 11).Value = 5929.1665
$ws.Cells.Item(96, 12).Value = 5869
$ws.Cells.Item(96, 13).Value = -4556.1665
$ws.Cells.Item(96, 14).Value = -8615
$ws.Cells.Item(122, 8).Value = 10874077
$ws.Cells.Item(122, 9).Value = 5231.875
$ws.Cells.Item(122, 10).Value = 35717150
$ws.Cells.Item(122, 11).Value = 15695.625
$ws.Cells.Item(122, 12).Value = 107151450
$ws.Cells.Item(122, 13).Value = -13245.625
$ws.Cells.Item(122, 14).Value = -107156350
$ws.Cells.Item(126, 8).Value = 2899
$ws.Cells.Item(126, 9).Value = 2899
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8697
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -6227
$ws.Cells.Item(126, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 9806960
$ws.Cells.Item(132, 9).Value = 15154707
$ws.Cells.Item(132, 10).Value = 2755.6667
$ws.Cells.Item(132, 11).Value = 45464121
$ws.Cells.Item(132, 12).Value = 8267.000100000001
$ws.Cells.Item(132, 13).Value = -45461591
$ws.Cells.Item(132, 14).Value = -13327.0001
$ws.Cells.Item(136, 8).Value = 22729384
$ws.Cells.Item(136, 9).Value = 25641966
$ws.Cells.Item(136, 11).Value = 76925898
$ws.Cells.Item(136, 13).Value = -76923348
